# Auto-generated edit script: updates crypto price/volume table cells
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to store $value as literal text even when it
    # parses as a number (Excel would otherwise coerce "202.96" etc.
    # into a numeric cell). A leading apostrophe forces text entry;
    # resetting the Style back to Normal afterwards keeps the cell
    # format identical to its un-touched neighbours.
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

$ws.Range('D2').Value = '70.775.97'
$ws.Range('E2').Value = '  +1.31%  '
$ws.Range('D3').Value = '3.610.21'
$ws.Range('E3').Value = '  +2.26%  '
$ws.Range('E4').Value = '  +0.05%  '
Set-TextValue $ws.Range('D5') '202.96'
$ws.Range('E5').Value = '  +3.71%  '
Set-TextValue $ws.Range('D6') '601.57'
$ws.Range('E6').Value = '  -0.74%  '
$ws.Range('E7').Value = '  +0.69%  '
Set-TextValue $ws.Range('D8') '0.999'
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('E9').Value = '  +6.46%  '
$ws.Range('E10').Value = '  +0.16%  '
Set-TextValue $ws.Range('D11') '54.12'
$ws.Range('E11').Value = '  +1.32%  '
$ws.Range('E12').Value = '  +0.30%  '
Set-TextValue $ws.Range('D13') '9.63'
$ws.Range('E13').Value = '  +1.71%  '
$ws.Range('D14').Value = '4.181.29'
$ws.Range('E14').Value = '  +2.14%  '
Set-TextValue $ws.Range('D15') '681.97'
$ws.Range('E15').Value = '  +14.46%  '
$ws.Range('D16').Value = '70.834.17'
$ws.Range('E16').Value = '  +1.27%  '
Set-TextValue $ws.Range('D17') '12.86'
$ws.Range('E17').Value = '  +1.23%  '
Set-TextValue $ws.Range('D18') '19.19'
$ws.Range('E18').Value = '  +1.11%  '
$ws.Range('D19').Value = '3.606.72'
$ws.Range('E19').Value = '  +2.62%  '
$ws.Range('E20').Value = '  +0.55%  '
$ws.Range('E21').Value = '  +1.68%  '
Set-TextValue $ws.Range('D22') '18.93'
$ws.Range('E22').Value = '  +5.77%  '
Set-TextValue $ws.Range('D23') '110.09'
$ws.Range('E23').Value = '  +6.62%  '
$ws.Range('E24').Value = '  +4.19%  '
$ws.Range('E25').Value = '  +0.38%  '
Set-TextValue $ws.Range('D26') '3.05'
$ws.Range('E26').Value = '  +0.04%  '
Set-TextValue $ws.Range('D27') '10.62'
$ws.Range('E27').Value = '  -1.52%  '
$ws.Range('E28').Value = '  -0.68%  '
Set-TextValue $ws.Range('D29') '10.14'
$ws.Range('E29').Value = '  +6.56%  '
Set-TextValue $ws.Range('D30') '34.45'
$ws.Range('E30').Value = '  +3.56%  '
Set-TextValue $ws.Range('D31') '4.48'
$ws.Range('E31').Value = '  +5.97%  '
$ws.Range('E32').Value = '  +1.72%  '
Set-TextValue $ws.Range('D33') '12.31'
$ws.Range('E33').Value = '  -0.14%  '
$ws.Range('E34').Value = '  +0.15%  '
$ws.Range('E35').Value = '  +0.21%  '
$ws.Range('D36').Value = '0.0₃0857'
$ws.Range('E36').Value = '  +5.57%  '
$ws.Range('D37').Value = '3.886.59'
$ws.Range('E37').Value = '  +2.79%  '
$ws.Range('E38').Value = '  -0.09%  '
Set-TextValue $ws.Range('D39') '514.10'
$ws.Range('E39').Value = '  +0.79%  '
Set-TextValue $ws.Range('D40') '3.03'
$ws.Range('E40').Value = '  -4.40%  '
$ws.Range('E41').Value = '  +1.04%  '
Set-TextValue $ws.Range('D42') '36.99'
$ws.Range('E42').Value = '  +1.64%  '
$ws.Range('E43').Value = '  +5.06%  '
Set-TextValue $ws.Range('D44') '0.385'
$ws.Range('E44').Value = '  -1.55%  '
Set-TextValue $ws.Range('D45') '0.0469'
$ws.Range('E45').Value = '  +4.37%  '
Set-TextValue $ws.Range('D46') '3.07'
$ws.Range('E46').Value = '  +8.94%  '
Set-TextValue $ws.Range('D47') '3.40'
$ws.Range('E47').Value = '  +4.84%  '
$ws.Range('E48').Value = '  +2.03%  '
Set-TextValue $ws.Range('D49') '8.65'
$ws.Range('E49').Value = '  +2.14%  '
$ws.Range('E50').Value = '  -0.24%  '
$ws.Range('B51').Value = 'CoreDAO'
$ws.Range('C51').Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
Set-TextValue $ws.Range('D51') '2.79'
$ws.Range('E51').Value = '  +68.12%  '
